$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "Institution N" -> "Author N Institution" ---
$ws.Range("E1").Value = "Author 1 Institution"
$ws.Range("G1").Value = "Author 2 Institution"
$ws.Range("I1").Value = "Author 3 Institution"
$ws.Range("K1").Value = "Author 4 Institution"
$ws.Range("M1").Value = "Author 5 Institution"
$ws.Range("O1").Value = "Author 6 Institution"
$ws.Range("Q1").Value = "Author 7 Institution"
$ws.Range("S1").Value = "Author 8 Institution"
$ws.Range("U1").Value = "Author 9 Institution"

# --- Column A (ID) on rows 2-63: change from numeric 0 to text ID values ---
$idMap = @{
    2 = "137458"
    3 = "137459"
    4 = "137460"
    5 = "137464"
    6 = "137465"
    7 = "137466"
    8 = "137468"
    9 = "137469"
    10 = "137470"
    11 = "137471"
    12 = "137472"
    13 = "137473"
    14 = "137474"
    15 = "137476"
    16 = "137477"
    17 = "137478"
    18 = "137479"
    19 = "137480"
    20 = "137486"
    21 = "137487"
    22 = "137555"
    23 = "137556"
    24 = "139165"
    25 = "139166"
    26 = "139167"
    27 = "139168"
    28 = "139169"
    29 = "139170"
    30 = "139172"
    31 = "139173"
    32 = "139174"
    33 = "139175"
    34 = "139176"
    35 = "139177"
    36 = "139178"
    37 = "139180"
    38 = "139181"
    39 = "139183"
    40 = "139184"
    41 = "139185"
    42 = "139186"
    43 = "139187"
    44 = "139188"
    45 = "139190"
    46 = "139192"
    47 = "139193"
    48 = "139194"
    49 = "139195"
    50 = "139189"
    51 = "139191"
    52 = "139179"
    53 = "139171"
    54 = "137467"
    55 = "139182"
    56 = "137461"
    57 = "137462"
    58 = "137475"
    59 = "139532"
    60 = "139533"
    61 = "139534"
    62 = "139531"
    63 = "137463"
}

foreach ($r in $idMap.Keys) {
    $cell = $ws.Range("A$r")
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $idMap[$r]
    $cell.Style = $origStyle
}
